# fix 27/4/2024 lần 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values (order chosen to match the final shared-string table order) ---
$ws.Range("A2").Value = "STT"
$ws.Range("A1").Value = "BÁO CÁO HÀNG TỒN KHO [time]"
$ws.Range("D2").Value = "Màu sắc"
$ws.Range("B2").Value = "Mã sản phẩm"
$ws.Range("C2").Value = "Tên sản phẩm"
$ws.Range("E2").Value = "Kích cỡ"
$ws.Range("F2").Value = "Số lượng hiện có "

# --- Header row (row 2) formatting: bold 12pt Times New Roman, thin box border, centered ---
$hdr = $ws.Range("A2:F2")
$hdr.Font.Name = "Times New Roman"
$hdr.Font.Bold = $true
$hdr.Font.Size = 12
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108

# --- Merge the title row first so every underlying cell shares one uniform style ---
$ttl = $ws.Range("A1:F1")
$ttl.Merge()

# --- Title row (row 1) formatting: bold 18pt Times New Roman, thin box border, centered ---
$ttl.Font.Name = "Times New Roman"
$ttl.Font.Bold = $true
$ttl.Font.Size = 18
$ttl.Borders.LineStyle = 1
$ttl.HorizontalAlignment = -4108
$ttl.VerticalAlignment = -4108

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 22.8
$ws.Rows.Item(2).RowHeight = 15.6

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 10.33203125
$ws.Columns.Item(2).ColumnWidth = 18.5546875
$ws.Columns.Item(3).ColumnWidth = 42.88671875
$ws.Columns.Item(4).ColumnWidth = 21.33203125
$ws.Columns.Item(5).ColumnWidth = 25.5546875
$ws.Columns.Item(6).ColumnWidth = 26.5546875

# --- Selection ---
$ws.Range("K7").Select()
